$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 11425
$ws1.Range("F19").Value = 11285
$ws1.Range("F20").Value = 11171
$ws1.Range("F25").Value = 40

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 11425
$ws4.Range("F20").Value = 11285
$ws4.Range("F21").Value = 11171
$ws4.Range("F26").Value = 40
